# Rebuild the LR-pair table (Sheet1, rows 2-17) with the updated natmi
# output that adds the "M2" sending-cluster rows (per Dr Hou's advice).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
  @("ECs", "Cx3cl1", "Itgav", "ECs", 3, 1, 9.355195333333333, 28.065586, 0.4481020953400846, 0.4481020953400845, 3, 1, 13.441269, 40.323807, 0.0897308213348123, 0.08973082133481232, 125.745697022878, 1131.711273205902, 0.04020856905671616, 0.04020856905671616),
  @("ECs", "Cx3cl1", "Itgav", "FAPs", 3, 1, 9.355195333333333, 28.065586, 0.4481020953400846, 0.4481020953400845, 3, 1, 54.711535, 164.134605, 0.3652416280068742, 0.3652416280068742, 511.8370969115033, 4606.533872203529, 0.163665538815304, 0.163665538815304),
  @("ECs", "Cx3cl1", "Itgav", "M2", 3, 1, 9.355195333333333, 28.065586, 0.4481020953400846, 0.4481020953400845, 3, 1, 63.67711, 191.03133, 0.4250937452800914, 0.4250937452800915, 595.7118023121533, 5361.40622080938, 0.1904853979759732, 0.1904853979759732),
  @("ECs", "Cx3cl1", "Itgav", "sCs", 3, 1, 9.355195333333333, 28.065586, 0.4481020953400846, 0.4481020953400845, 3, 1, 17.96553866666667, 53.896616, 0.119933805378222, 0.119933805378222, 168.0711234952196, 1512.640111456976, 0.05374258949209119, 0.05374258949209119),
  @("FAPs", "Cx3cl1", "Itgav", "ECs", 3, 1, 10.423773, 31.271319, 0.4992856221832745, 0.4992856221832744, 3, 1, 13.441269, 40.323807, 0.0897308213348123, 0.08973082133481232, 140.108736887937, 1260.978631991433, 0.044801308959168, 0.04480130895916801),
  @("FAPs", "Cx3cl1", "Itgav", "FAPs", 3, 1, 10.423773, 31.271319, 0.4992856221832745, 0.4992856221832744, 3, 1, 54.711535, 164.134605, 0.3652416280068742, 0.3652416280068742, 570.3006213215549, 5132.705591893994, 0.1823598934866442, 0.1823598934866443),
  @("FAPs", "Cx3cl1", "Itgav", "M2", 3, 1, 10.423773, 31.271319, 0.4992856221832745, 0.4992856221832744, 3, 1, 63.67711, 191.03133, 0.4250937452800914, 0.4250937452800915, 663.75573993603, 5973.80165942427, 0.2122431950983888, 0.2122431950983888),
  @("FAPs", "Cx3cl1", "Itgav", "sCs", 3, 1, 10.423773, 31.271319, 0.4992856221832745, 0.4992856221832744, 3, 1, 17.96553866666667, 53.896616, 0.119933805378222, 0.119933805378222, 187.268696884056, 1685.418271956504, 0.05988122463907333, 0.05988122463907333),
  @("M2", "Cx3cl1", "Itgav", "ECs", 1, 0.3333333333333333, 0.009880666666666668, 0.029642, 0.0004732715115968286, 0.0004732715115968285, 3, 1, 13.441269, 40.323807, 0.0897308213348123, 0.08973082133481232, 0.132808698566, 1.195278287094, 0.00004246704144995158, 0.00004246704144995157),
  @("M2", "Cx3cl1", "Itgav", "FAPs", 1, 0.3333333333333333, 0.009880666666666668, 0.029642, 0.0004732715115968286, 0.0004732715115968285, 3, 1, 54.711535, 164.134605, 0.3652416280068742, 0.3652416280068742, 0.5405864401566667, 4.86527796141, 0.0001728584573848999, 0.0001728584573848999),
  @("M2", "Cx3cl1", "Itgav", "M2", 1, 0.3333333333333333, 0.009880666666666668, 0.029642, 0.0004732715115968286, 0.0004732715115968285, 3, 1, 63.67711, 191.03133, 0.4250937452800914, 0.4250937452800915, 0.6291722982066668, 5.66255068386, 0.0002011847593990661, 0.0002011847593990661),
  @("M2", "Cx3cl1", "Itgav", "sCs", 1, 0.3333333333333333, 0.009880666666666668, 0.029642, 0.0004732715115968286, 0.0004732715115968285, 3, 1, 17.96553866666667, 53.896616, 0.119933805378222, 0.119933805378222, 0.1775114990524445, 1.597603491472, 0.00005676125336291097, 0.00005676125336291097),
  @("sCs", "Cx3cl1", "Itgav", "ECs", 3, 1, 1.088525666666667, 3.265577, 0.05213901096504407, 0.05213901096504406, 3, 1, 13.441269, 40.323807, 0.0897308213348123, 0.08973082133481232, 14.631166299071, 131.680496691639, 0.00467847627747819, 0.004678476277478189),
  @("sCs", "Cx3cl1", "Itgav", "FAPs", 3, 1, 1.088525666666667, 3.265577, 0.05213901096504407, 0.05213901096504406, 3, 1, 54.711535, 164.134605, 0.3652416280068742, 0.3652416280068742, 59.55491011023165, 535.9941909920849, 0.01904333724754096, 0.01904333724754096),
  @("sCs", "Cx3cl1", "Itgav", "M2", 3, 1, 1.088525666666667, 3.265577, 0.05213901096504407, 0.05213901096504406, 3, 1, 63.67711, 191.03133, 0.4250937452800914, 0.4250937452800915, 69.31416861415666, 623.8275175274099, 0.02216396744633034, 0.02216396744633033),
  @("sCs", "Cx3cl1", "Itgav", "sCs", 3, 1, 1.088525666666667, 3.265577, 0.05213901096504407, 0.05213901096504406, 3, 1, 17.96553866666667, 53.896616, 0.119933805378222, 0.119933805378222, 19.55594995415911, 176.003549587432, 0.006253229993694579, 0.006253229993694578)
)

$startRow = 2
for ($i = 0; $i -lt $rows.Length; $i++) {
  $r = $startRow + $i
  $rowValues = $rows[$i]
  for ($c = 0; $c -lt $rowValues.Length; $c++) {
    $ws.Cells.Item($r, $c + 1).Value = $rowValues[$c]
  }
}
